# Weekly update: insert a new price record for "Arveja Verde" at row 258,
# shifting the existing rows 258-283 down to 259-284 (the last existing
# row, old row 283, becomes the new row 284).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 258 - this pushes rows 258..283 down to
# 259..284 and leaves a blank row 258 ready to be populated, mirroring the
# Excel UI action of right-clicking a row header and choosing "Insert".
$ws.Rows.Item(258).Insert()

# Populate the newly inserted row 258 with this week's record.
$ws.Cells.Item(258, 1).Value = 6
$ws.Cells.Item(258, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(258, 3).Value = "Metropolitana"
$ws.Cells.Item(258, 4).Value = 44918
$ws.Cells.Item(258, 5).Value = 13
$ws.Cells.Item(258, 6).Value = 100112022
$ws.Cells.Item(258, 7).Value = "Arveja Verde"
$ws.Cells.Item(258, 8).Value = "Perfection"
$ws.Cells.Item(258, 9).Value = "Primera"
$ws.Cells.Item(258, 10).Value = 350
$ws.Cells.Item(258, 11).Value = 20000
$ws.Cells.Item(258, 12).Value = 22000
$ws.Cells.Item(258, 13).Value = 21143
$ws.Cells.Item(258, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(258, 15).Value = "Carahue"
$ws.Cells.Item(258, 16).Value = 846
$ws.Cells.Item(258, 17).Value = 25
$ws.Cells.Item(258, 18).Value = "Hortaliza"
